# practica 3 finished, subject nearly passed
# Extends the FSM state table (T1/T2) with a new state "101 (F)" and mirrors
# it, together with the associated 4-bit encoded transition table, into a
# new block of columns (AO:AV) on Sheet1. Also enlarges/repositions the
# explanatory picture and updates the view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New T1 / T2 tables (state name / next-state / output) in AO:AV,
#    rows 22-29, replicating the existing T1 (A17:C23) and T2 (F17:H23)
#    tables but for the 6-state machine that includes "101 (F)".
# ---------------------------------------------------------------------

# Row 22: section headers "T1" / "T2"
$ws.Range("A17").Copy()
$ws.Range("AO22").PasteSpecial(-4122)
$ws.Range("AO22").Value = "T1"

$ws.Range("F17").Copy()
$ws.Range("AT22").PasteSpecial(-4122)
$ws.Range("AT22").Value = "T2"

# Row 23: column headers
$ws.Range("A18:D18").Copy()
$ws.Range("AO23:AR23").PasteSpecial(-4122)
$ws.Range("AO23").Value = "Z(n) \ X (n)"
$ws.Range("AP23").Value = "0"
$ws.Range("AQ23").Value = "1"

$ws.Range("F18:H18").Copy()
$ws.Range("AT23:AV23").PasteSpecial(-4122)
$ws.Range("AT23").Value = "Z(n) \ X (n)"
$ws.Range("AU23").Value = "Y(n)"

# Rows 24-28: states A-E, copied from the existing T1/T2 body rows 19-23
$ws.Range("A19:D19").Copy()
$ws.Range("AO24:AR24").PasteSpecial(-4122)
$ws.Range("AO24").Value = "000 ( A )"
$ws.Range("AP24").Value = "000"
$ws.Range("AQ24").Value = "001"

$ws.Range("F19:H19").Copy()
$ws.Range("AT24:AV24").PasteSpecial(-4122)
$ws.Range("AT24").Value = "000 ( A )"
$ws.Range("AU24").Value = "0"

$ws.Range("A20:D20").Copy()
$ws.Range("AO25:AR25").PasteSpecial(-4122)
$ws.Range("AO25").Value = "001 ( B )"
$ws.Range("AP25").Value = "010"
$ws.Range("AQ25").Value = "001"

$ws.Range("F20:H20").Copy()
$ws.Range("AT25:AV25").PasteSpecial(-4122)
$ws.Range("AT25").Value = "001 ( B )"
$ws.Range("AU25").Value = "0"

$ws.Range("A21:D21").Copy()
$ws.Range("AO26:AR26").PasteSpecial(-4122)
$ws.Range("AO26").Value = "010 ( C )"
$ws.Range("AP26").Value = "000"
$ws.Range("AQ26").Value = "011"

$ws.Range("F21:H21").Copy()
$ws.Range("AT26:AV26").PasteSpecial(-4122)
$ws.Range("AT26").Value = "010 ( C )"
$ws.Range("AU26").Value = "0"

$ws.Range("A22:D22").Copy()
$ws.Range("AO27:AR27").PasteSpecial(-4122)
$ws.Range("AO27").Value = "011 ( D )"
$ws.Range("AP27").Value = "100"
$ws.Range("AQ27").Value = "001"

$ws.Range("F22:H22").Copy()
$ws.Range("AT27:AV27").PasteSpecial(-4122)
$ws.Range("AT27").Value = "011 ( D )"
$ws.Range("AU27").Value = "0"

$ws.Range("A23:D23").Copy()
$ws.Range("AO28:AR28").PasteSpecial(-4122)
$ws.Range("AO28").Value = "100 ( E ) "
$ws.Range("AP28").Value = "101"
$ws.Range("AQ28").Value = "001"

$ws.Range("F23:H23").Copy()
$ws.Range("AT28:AV28").PasteSpecial(-4122)
$ws.Range("AT28").Value = "100 ( E ) "
$ws.Range("AU28").Value = "0"

# Row 29: new state "101 (F)" - distinct formatting (centered box, thin
# left/right border only) since it is the newly added state.
$ws.Range("AO29").Value = "101 (F)"
$ws.Range("AO29").NumberFormat = "@"
$ws.Range("AO29").HorizontalAlignment = -4108
$ws.Range("AO29").VerticalAlignment = -4108
$ws.Range("AO29").Borders.Item(7).LineStyle = 1
$ws.Range("AO29").Borders.Item(7).Weight = 2
$ws.Range("AO29").Borders.Item(7).ColorIndex = 1
$ws.Range("AO29").Borders.Item(10).LineStyle = 1
$ws.Range("AO29").Borders.Item(10).Weight = 2
$ws.Range("AO29").Borders.Item(10).ColorIndex = 1

$ws.Range("A8:B8").Copy()
$ws.Range("AP29:AQ29").PasteSpecial(-4122)
$ws.Range("AP29").Value = "000"
$ws.Range("AQ29").Value = "001"

$ws.Range("AO29").Copy()
$ws.Range("AT29").PasteSpecial(-4122)
$ws.Range("AT29").Value = "101 (F)"
$ws.Range("AU29").Value = 1

# ---------------------------------------------------------------------
# 2) New 4-bit encoded transition table, rows 33-38 (mirrors C40:D44
#    which already exists a few rows below), columns AP:AQ.
# ---------------------------------------------------------------------

$ws.Range("C40:D40").Copy()
$ws.Range("AP33:AQ33").PasteSpecial(-4122)
$ws.Range("AP33").Value = "0001"
$ws.Range("AQ33").Value = "0000"

$ws.Range("C41:D41").Copy()
$ws.Range("AP34:AQ34").PasteSpecial(-4122)
$ws.Range("AP34").Value = "0001"
$ws.Range("AQ34").Value = "0010"

$ws.Range("C42:D42").Copy()
$ws.Range("AP35:AQ35").PasteSpecial(-4122)
$ws.Range("AP35").Value = "0011"
$ws.Range("AQ35").Value = "0000"

$ws.Range("C43:D43").Copy()
$ws.Range("AP36:AQ36").PasteSpecial(-4122)
$ws.Range("AP36").Value = "0001"
$ws.Range("AQ36").Value = "0100"

$ws.Range("C44:D44").Copy()
$ws.Range("AP37:AQ37").PasteSpecial(-4122)
$ws.Range("AP37").Value = "0001"
$ws.Range("AQ37").Value = "0101"

$ws.Range("C46:D46").Copy()
$ws.Range("AP38:AQ38").PasteSpecial(-4122)
$ws.Range("AP38").Value = "0001"
$ws.Range("AQ38").Value = "0000"

# ---------------------------------------------------------------------
# 3) Resize / reposition the explanatory picture (3rd picture on the
#    sheet) so it also covers the newly added table.
# ---------------------------------------------------------------------

$pic = $ws.Shapes.Item(3)
$pic.Left = 13009417 / 9525
$pic.Top = 2687780 / 9525
$pic.Width = 11838015 / 9525
$pic.Height = 6719455 / 9525

# ---------------------------------------------------------------------
# 4) View state: selected cell and zoom/scroll position.
# ---------------------------------------------------------------------

$ws.Range("AP29:AQ29").Select()
$excel.ActiveWindow.Zoom = 160
$excel.ActiveWindow.ScrollColumn = 38
$excel.ActiveWindow.ScrollRow = 22
